$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "72.208.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.82%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "4.030.88"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.17%  "

# Row 4
$ws.Range("E4").Value = "  +0.27%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.13%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.64"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.99%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.705"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +13.33%  "

# Row 8
$ws.Range("E8").Value = "  +0.14%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.755"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.10%  "

# Row 10
$ws.Range("E10").Value = "  -1.09%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000326"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.02%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.63"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +11.40%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.74"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.22%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.675.49"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.93%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.020.19"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.07%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.90%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.58"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.73%  "

# Row 18
$ws.Range("E18").Value = "  -0.57%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.19%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.138.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.04%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "430.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.53%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "99.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +10.48%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.09%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.84%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.52"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.86%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.11"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.71%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.91"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.66%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +28.28%  "

# Row 29
$ws.Range("E29").Value = "  +2.06%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.10"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.71%  "

# Row 31
$ws.Range("E31").Value = "  +0.35%  "

# Row 32
$ws.Range("E32").Value = "  +2.86%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "681.15"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.17%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.97"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.02%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "66.32"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.33%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "42.88"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.79%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.427"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.81%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.154"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.00%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.51"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +11.38%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0824"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -10.18%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.42"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.36%  "

# Row 42
$ws.Range("E42").Value = "  -0.10%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.09%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0490"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.40%  "

# Row 45
$ws.Range("E45").Value = "  +6.33%  "

# Row 46
$ws.Range("B46").Value = "THORChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.65"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.14%  "

# Row 47
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.39"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.73%  "

# Row 48
$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -9.24%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.02"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.78%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.38"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.31%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "146.41"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.87%  "
